$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry (row 17): date, task description, hours spent.
$ws.Range("A17").Value = " `r`n18.07.2019"
$ws.Range("B17").Value = "Refractored the code a bit. Added a hidden button on the card images to flip them over in case they are double-sided. The button is visible only when mouse enters a card which is double-sided. Clicking the button switches to the backside"
$ws.Range("C17").Value = 5

# Match the wrap-text formatting used by the rest of the table.
$ws.Range("A17:B17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 60

# Scroll/selection state, as left by the author after adding the row.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D17").Select()
